# Edit workbook: rename "Troof" -> "Troof1", and append six new sheets
# (T6, T7, T8, T9, T10, Troof2) after "Twall", each populated with the
# data tables described by the upstream diff.

$wb = $excel.ActiveWorkbook

# 1. Rename the existing "Troof" sheet to "Troof1".
$wb.Worksheets("Troof").Name = "Troof1"

# xlPasteFormats - used to copy just the header cell-style (bold font,
# thin border, center/top alignment) from the matching "template" sheet
# onto the new sheet's header row, so the new sheets look consistent
# with the rest of the workbook.
$xlPasteFormats = -4122
# Excel pads ColumnWidth by this much versus the width value stored in
# the saved file, so subtract it to land on the exact target width.
$widthPad = 0.8333333333333333

function Set-ColWidths {
    param($ws, [double[]]$widths)
    for ($i = 0; $i -lt $widths.Length; $i++) {
        $ws.Columns.Item($i + 1).ColumnWidth = $widths[$i] - $widthPad
    }
}

function Fill-Sheet {
    param($ws, $data)
    for ($r = 0; $r -lt $data.Length; $r++) {
        $row = $data[$r]
        for ($c = 0; $c -lt $row.Length; $c++) {
            $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
        }
    }
}

# NB: deliberately takes no parameters - it reads $script:_ws /
# $script:_templateWs / $script:_colCount instead. Passing both the
# source (template) and destination worksheet COM objects into the same
# function as parameters causes the Copy()/PasteSpecial() clipboard
# hand-off to silently lose the pasted style; module-scoped variables
# set right before the call avoid that.
function Style-Header {
    $script:_templateWs.Range($script:_templateWs.Cells.Item(1,1), $script:_templateWs.Cells.Item(1, $script:_colCount)).Copy() | Out-Null
    $script:_ws.Range($script:_ws.Cells.Item(1,1), $script:_ws.Cells.Item(1, $script:_colCount)).PasteSpecial($xlPasteFormats) | Out-Null
}

$prevSheet = $wb.Worksheets("Twall")

# ---- T6 (mirrors T1: b_m, d_m, h_m, e_m, ze_m, crze) ----
$t6 = $wb.Worksheets.Add($null, $prevSheet)
$t6.Name = "T6"
$t6Data = @(
    @('b_m', 'd_m', 'h_m', 'e_m', 'ze_m', 'crze'),
    @(18, 16, 9, 18, 9, 0.7312574370573635)
)
Fill-Sheet $t6 $t6Data
$script:_ws = $t6
$script:_templateWs = $wb.Worksheets("T1")
$script:_colCount = 6
Style-Header
Set-ColWidths $t6 @(5, 5, 5, 5, 6, 20)
$prevSheet = $t6

# ---- T7 (mirrors T2: Ivze, Ceze, qpze_N_m_2) ----
$t7 = $wb.Worksheets.Add($null, $prevSheet)
$t7.Name = "T7"
$t7Data = @(
    @('Ivze', 'Ceze', 'qpze_N_m_2'),
    @(0.294014103795206, 1.635279882023036, 817.6399410115181)
)
Fill-Sheet $t7 $t7Data
$script:_ws = $t7
$script:_templateWs = $wb.Worksheets("T2")
$script:_colCount = 3
Style-Header
Set-ColWidths $t7 @(19, 19, 19)
$prevSheet = $t7

# ---- T8 (mirrors T3: c, epf, ipf1, we, wi1) ----
$t8 = $wb.Worksheets.Add($null, $prevSheet)
$t8.Name = "T8"
$t8Data = @(
    @('c', 'epf', 'ipf1', 'we', 'wi1'),
    @('D', 0.8, -0.1077083333333334, 654.1119528092145, -88.06663531311565),
    @('E', -0.3, -0.1077083333333334, -245.2919823034554, -88.06663531311565),
    @('A', -0.3, -0.1077083333333334, -245.2919823034554, -88.06663531311565),
    @('B', -1, -0.1077083333333334, -817.6399410115181, -88.06663531311565),
    @('F1', -1.42886653550543, -0.1077083333333334, -1168.298349803992, -88.06663531311565),
    @('G1', -1.3, -0.1077083333333334, -1062.931923314974, -88.06663531311565),
    @('G2', -1.3, -0.1077083333333334, -1062.931923314974, -88.06663531311565),
    @('F2', -1.42886653550543, -0.1077083333333334, -1168.298349803992, -88.06663531311565),
    @('H1', -0.6247119229084849, -0.1077083333333334, -510.7894197960856, -88.06663531311565),
    @('H2', -0.6247119229084849, -0.1077083333333334, -510.7894197960856, -88.06663531311565),
    @('I1', -0.5247119229084849, -0.1077083333333334, -429.0254256949338, -88.06663531311565),
    @('I2', -0.5247119229084849, -0.1077083333333334, -429.0254256949338, -88.06663531311565)
)
Fill-Sheet $t8 $t8Data
$script:_ws = $t8
$script:_templateWs = $wb.Worksheets("T3")
$script:_colCount = 5
Style-Header
Set-ColWidths $t8 @(4, 21, 21, 20, 20)
$prevSheet = $t8

# ---- T9 (mirrors T4: c, Fwehor, Fwi1hor, Fwez, Fwi1z, x, y, z) ----
$t9 = $wb.Worksheets.Add($null, $prevSheet)
$t9.Name = "T9"
$t9Data = @(
    @('c', 'Fwehor', 'Fwi1hor', 'Fwez', 'Fwi1z', 'x', 'y', 'z'),
    @('D', 105966.1363550927, 0, 0, 0, 0, 9, 9),
    @('E', 39737.30113315978, -14266.79492072474, 0, 0, 16, 9, 9),
    @('F1', 0, 0, 9694.060717859842, -730.7408335261226, 2.25, 0.9, 9),
    @('G1', 0, 0, 17639.54661974082, -1461.481667052245, 2.25, 0.9, 9),
    @('G2', 0, 0, 17639.54661974082, -1461.481667052245, 11.25, 0.9, 9),
    @('F2', 0, 0, 9694.060717859842, -730.7408335261226, 15.75, 0.9, 9),
    @('H1', 0, 0, 58394.64747514817, -10067.98481747102, 4.5, 5.4, 9),
    @('H2', 0, 0, 58394.64747514817, -10067.98481747102, 13.5, 5.4, 9),
    @('I1', 0, 0, 49047.19542024277, -10067.98481747102, 4.5, 12.5, 9),
    @('I2', 0, 0, 49047.19542024277, -10067.98481747102, 13.5, 12.5, 9)
)
Fill-Sheet $t9 $t9Data
$script:_ws = $t9
$script:_templateWs = $wb.Worksheets("T4")
$script:_colCount = 8
Style-Header
Set-ColWidths $t9 @(4, 19, 20, 19, 20, 7, 6, 3)
$prevSheet = $t9

# ---- T10 (mirrors T5: row, Rhor, ZRhor, Rz, d, horRz, Mrv) ----
$t10 = $wb.Worksheets.Add($null, $prevSheet)
$t10.Name = "T10"
$t10Data = @(
    @('row', 'Rhor', 'ZRhor', 'Rz', 'd', 'horRz', 'Mrv'),
    @('row', 131436.6425675278, 9, 224894.5161949423, 16, 6.855048049026271, 3282142.780811639)
)
Fill-Sheet $t10 $t10Data
$script:_ws = $t10
$script:_templateWs = $wb.Worksheets("T5")
$script:_colCount = 7
Style-Header
Set-ColWidths $t10 @(5, 19, 7, 19, 4, 19, 19)
$prevSheet = $t10

# ---- Troof2 (mirrors Troof1: factors, F, G, H, I) ----
$troof2 = $wb.Worksheets.Add($null, $prevSheet)
$troof2.Name = "Troof2"
$troof2Data = @(
    @('factors', 'F', 'G', 'H', 'I'),
    @('epf1-', -2.04942384581697, -2, -1.2, -0.5247119229084849),
    @('epf10-', -1.374135768725455, -1.3, -0.6247119229084849, -0.5247119229084849)
)
Fill-Sheet $troof2 $troof2Data
$script:_ws = $troof2
$script:_templateWs = $wb.Worksheets("Troof1")
$script:_colCount = 5
Style-Header
Set-ColWidths $troof2 @(9, 20, 6, 21, 21)
$prevSheet = $troof2

Write-Host "Edit complete"
